# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) previously held a different stat ("Strike#"); this
# regenerates it to hold true strikeout counts (K) per game, row by row.
# Row 43 additionally needed its IP (H) and I0 (I) values corrected as part
# of the same regeneration pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    8  = 2
    9  = 2
    10 = 1
    11 = 0
    12 = 0
    13 = 1
    14 = 1
    15 = 2
    16 = 1
    17 = 1
    18 = 1
    19 = 2
    20 = 2
    21 = 1
    22 = 0
    23 = 0
    24 = 1
    25 = 2
    26 = 1
    27 = 2
    28 = 1
    29 = 0
    30 = 2
    31 = 0
    32 = 2
    33 = 0
    34 = 0
    35 = 2
    36 = 0
    37 = 1
    38 = 1
    39 = 1
    40 = 1
    41 = 0
    42 = 1
    43 = 0
    44 = 1
    45 = 3
    46 = 1
    47 = 3
    48 = 2
    49 = 2
    50 = 1
    51 = 0
    52 = 1
    53 = 2
    54 = 1
    55 = 1
    56 = 1
    57 = 0
    59 = 1
    60 = 0
    61 = 2
    62 = 1
    63 = 1
    64 = 1
    65 = 1
    66 = 0
    67 = 0
    68 = 0
    69 = 0
    70 = 1
    71 = 1
    72 = 0
    73 = 1
    74 = 1
    75 = 2
    76 = 2
    77 = 1
    78 = 3
    79 = 2
    80 = 2
    81 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}

# Row 43 also needed IP (H) and I0 (I) corrected during the K regeneration.
$ws.Range("H43").Value = 2
$ws.Range("I43").Value = 8
